# Update edited session - 2025-11-20T11:55:48.292Z - Cache Bust ID: 1763639748292auif45tys
#
# 1) Rename the worksheet from "Physiology" to "Session".
# 2) Remove the log row for student 231521 (originally row 17) - the
#    scanner record that preceded 231479 - by deleting the entire row,
#    which shifts every subsequent row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab.
$ws.Name = "Session"

# Delete the row for student 231521 (13:42:32 Scan entry) - everything
# below it (rows 18-24) shifts up to fill the gap.
$ws.Rows("17").Delete()
